$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last annotation" (row 2) with corrected position/size values
$ws.Range("C2").Value = "#caff2b0e-3932-429f-b2fb-fc53d3980941"
$ws.Range("D2").Value = 1213
$ws.Range("E2").Value = 152
$ws.Range("F2").Value = 192
$ws.Range("G2").Value = 36
$ws.Range("H2").Value = $false

# Remove the now-obsolete trailing annotation row entirely
$ws.Rows(3).Delete()
